$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.677.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.39%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.91%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6
$ws.Range("E6").Value = "  +2.73%  "

# Row 9
$ws.Range("E9").Value = "  +1.35%  "

# Row 10
$ws.Range("E10").Value = "  +1.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.91%  "

# Row 14
$ws.Range("E14").Value = "  +2.16%  "

# Row 15
$ws.Range("E15").Value = "  +1.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.688.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.64%  "

# Row 20
$ws.Range("E20").Value = "  +0.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.16%  "

# Row 23
$ws.Range("E23").Value = "  +2.53%  "

# Row 24
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.44%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("E28").Value = "  +4.50%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.44%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0511"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.68%  "

# Row 31
$ws.Range("E31").Value = "  +0.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.22%  "

# Row 33
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("E34").Value = "  +1.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.200.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "

# Row 37
$ws.Range("E37").Value = "  +5.92%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.812"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "

# Row 39
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.506"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.86%  "

# Row 41
$ws.Range("E41").Value = "  -0.56%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.41%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.795"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.74%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.774.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.55%  "

# Row 48
$ws.Range("E48").Value = "  +0.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.63%  "

# Row 50
$ws.Range("E50").Value = "  +0.50%  "

# Row 51
$ws.Range("E51").Value = "  +0.22%  "
